$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 from numeric 0 to text "2h 11m p3" (keeps existing style)
$ws.Range("B3").Value = "2h 11m p3"

# Move the active selection from C3 to B4
$ws.Range("B4").Select()
